# Updates cryptos list price/volume(1h) data to latest snapshot values.
# Two rows (29/30) also swap which coin occupies them (RenderToken <-> Binance-PegBSC-USD).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.107.88"
$ws.Range("E2").Value = "  -2.45%  "

$ws.Range("D3").Value = "3.513.70"
$ws.Range("E3").Value = "  -1.94%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.25"
$ws.Range("E5").Value = "  -2.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.21"
$ws.Range("E6").Value = "  -3.94%  "

$ws.Range("D7").Value = "3.514.00"
$ws.Range("E7").Value = "  -1.92%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -1.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.11"
$ws.Range("E11").Value = "  -1.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  -1.16%  "

$ws.Range("D13").Value = "4.087.48"
$ws.Range("E13").Value = "  -2.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.87"
$ws.Range("E14").Value = "  -0.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  -3.62%  "

$ws.Range("E16").Value = "  +0.03%  "

$ws.Range("D17").Value = "3.493.65"
$ws.Range("E17").Value = "  -2.54%  "

$ws.Range("D18").Value = "64.097.95"
$ws.Range("E18").Value = "  -2.60%  "

$ws.Range("E19").Value = "  +0.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.38"
$ws.Range("E20").Value = "  -1.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.67"
$ws.Range("E21").Value = "  -3.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.51"
$ws.Range("E22").Value = "  -1.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.579"
$ws.Range("E23").Value = "  -2.05%  "

$ws.Range("D24").Value = "3.644.15"
$ws.Range("E24").Value = "  -2.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.97"
$ws.Range("E25").Value = "  -1.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000113"
$ws.Range("E27").Value = "  -4.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.58"
$ws.Range("E28").Value = "  -3.70%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.46"
$ws.Range("E29").Value = "  -7.81%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.25"
$ws.Range("E31").Value = "  -5.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.23"
$ws.Range("E32").Value = "  -4.06%  "

$ws.Range("D33").Value = "3.507.38"
$ws.Range("E33").Value = "  -2.29%  "

$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.89"
$ws.Range("E35").Value = "  -2.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.145"
$ws.Range("E36").Value = "  -2.24%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.34"
$ws.Range("E37").Value = "  -1.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.97"
$ws.Range("E38").Value = "  -1.36%  "

$ws.Range("E39").Value = "  -3.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "168.74"
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0810"
$ws.Range("E41").Value = "  -3.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.814"
$ws.Range("E42").Value = "  -3.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.62"
$ws.Range("E43").Value = "  -1.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.96"
$ws.Range("E45").Value = "  -2.75%  "

$ws.Range("E46").Value = "  -5.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.38"
$ws.Range("E47").Value = "  -3.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.65"
$ws.Range("E48").Value = "  -3.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.90"
$ws.Range("E49").Value = "  -1.72%  "

$ws.Range("D50").Value = "2.444.43"
$ws.Range("E50").Value = "  -0.43%  "

$ws.Range("E51").Value = "  -0.97%  "
